# Fix userCount (column H) values that were mis-imported from J6 ("un par de cagadas al traer los datos de la J6")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 82
$ws.Range("H4").Value  = 271
$ws.Range("H5").Value  = 329
$ws.Range("H6").Value  = 86
$ws.Range("H7").Value  = 133
$ws.Range("H8").Value  = 93
$ws.Range("H9").Value  = 83
$ws.Range("H10").Value = 50
$ws.Range("H12").Value = 473
$ws.Range("H13").Value = 324
$ws.Range("H14").Value = 105
$ws.Range("H15").Value = 70
